# Update Pdgfb-Pdgfra NATMI LR-pair output with newly computed TPM-based values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 51.93629233333333
$ws.Range("H2").Value = 155.808877
$ws.Range("I2").Value = 0.7704232182162135
$ws.Range("J2").Value = 0.7704232182162134
$ws.Range("M2").Value = 0.667106
$ws.Range("N2").Value = 2.001318
$ws.Range("O2").Value = 0.003817114239487378
$ws.Range("P2").Value = 0.003817114239487378
$ws.Range("Q2").Value = 34.64701223332067
$ws.Range("R2").Value = 311.823110099886
$ws.Range("S2").Value = 0.0029407934366848
$ws.Range("T2").Value = 0.0029407934366848
$ws.Range("G3").Value = 51.93629233333333
$ws.Range("H3").Value = 155.808877
$ws.Range("I3").Value = 0.7704232182162135
$ws.Range("J3").Value = 0.7704232182162134
$ws.Range("O3").Value = 0.9945745510447523
$ws.Range("P3").Value = 0.9945745510447522
$ws.Range("Q3").Value = 9027.509913254427
$ws.Range("R3").Value = 81247.58921928985
$ws.Range("S3").Value = 0.7662433263718438
$ws.Range("T3").Value = 0.7662433263718437
$ws.Range("G4").Value = 51.93629233333333
$ws.Range("H4").Value = 155.808877
$ws.Range("I4").Value = 0.7704232182162135
$ws.Range("J4").Value = 0.7704232182162134
$ws.Range("K4").Value = 2.0
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.260372
$ws.Range("N4").Value = 0.7811159999999999
$ws.Range("O4").Value = 0.001489822709979835
$ws.Range("P4").Value = 0.001489822709979834
$ws.Range("Q4").Value = 13.52275630741467
$ws.Range("R4").Value = 121.704806766732
$ws.Range("S4").Value = 0.001147794006794265
$ws.Range("T4").Value = 0.001147794006794265
$ws.Range("G5").Value = 51.93629233333333
$ws.Range("H5").Value = 155.808877
$ws.Range("I5").Value = 0.7704232182162135
$ws.Range("J5").Value = 0.7704232182162134
$ws.Range("K5").Value = 2.0
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.020712
$ws.Range("N5").Value = 0.062136
$ws.Range("O5").Value = 0.0001185120057805845
$ws.Range("P5").Value = 0.0001185120057805844
$ws.Range("Q5").Value = 1.075704486808
$ws.Range("R5").Value = 9.681340381272
$ws.Range("S5").Value = 0.00009130440089073638
$ws.Range("T5").Value = 0.00009130440089073636
$ws.Range("G6").Value = 0.03171066666666666
$ws.Range("H6").Value = 0.095132
$ws.Range("I6").Value = 0.0004703961867034368
$ws.Range("J6").Value = 0.0004703961867034368
$ws.Range("M6").Value = 0.667106
$ws.Range("N6").Value = 2.001318
$ws.Range("O6").Value = 0.003817114239487378
$ws.Range("P6").Value = 0.003817114239487378
$ws.Range("Q6").Value = 0.02115437599733333
$ws.Range("R6").Value = 0.190389383976
$ws.Range("S6").Value = 0.000001795555982466252
$ws.Range("T6").Value = 0.000001795555982466252
$ws.Range("G7").Value = 0.03171066666666666
$ws.Range("H7").Value = 0.095132
$ws.Range("I7").Value = 0.0004703961867034368
$ws.Range("J7").Value = 0.0004703961867034368
$ws.Range("O7").Value = 0.9945745510447523
$ws.Range("P7").Value = 0.9945745510447522
$ws.Range("Q7").Value = 5.511913631645776
$ws.Range("R7").Value = 49.60722268481199
$ws.Range("S7").Value = 0.0004678440762037341
$ws.Range("T7").Value = 0.0004678440762037341
$ws.Range("G8").Value = 0.03171066666666666
$ws.Range("H8").Value = 0.095132
$ws.Range("I8").Value = 0.0004703961867034368
$ws.Range("J8").Value = 0.0004703961867034368
$ws.Range("K8").Value = 2.0
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.260372
$ws.Range("N8").Value = 0.7811159999999999
$ws.Range("O8").Value = 0.001489822709979835
$ws.Range("P8").Value = 0.001489822709979834
$ws.Range("Q8").Value = 0.008256569701333332
$ws.Range("R8").Value = 0.07430912731199998
$ws.Range("S8").Value = 0.0000007008069216386944
$ws.Range("T8").Value = 0.0000007008069216386943
$ws.Range("G9").Value = 0.03171066666666666
$ws.Range("H9").Value = 0.095132
$ws.Range("I9").Value = 0.0004703961867034368
$ws.Range("J9").Value = 0.0004703961867034368
$ws.Range("K9").Value = 2.0
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.020712
$ws.Range("N9").Value = 0.062136
$ws.Range("O9").Value = 0.0001185120057805845
$ws.Range("P9").Value = 0.0001185120057805844
$ws.Range("Q9").Value = 0.0006567913279999999
$ws.Range("R9").Value = 0.005911121952
$ws.Range("S9").Value = 0.00000005574759559776258
$ws.Range("T9").Value = 0.00000005574759559776258
$ws.Range("G10").Value = 2.162051666666667
$ws.Range("H10").Value = 6.486155
$ws.Range("I10").Value = 0.03207188515291837
$ws.Range("J10").Value = 0.03207188515291837
$ws.Range("M10").Value = 0.667106
$ws.Range("N10").Value = 2.001318
$ws.Range("O10").Value = 0.003817114239487378
$ws.Range("P10").Value = 0.003817114239487378
$ws.Range("Q10").Value = 1.442317639143333
$ws.Range("R10").Value = 12.98085875229
$ws.Range("S10").Value = 0.0001224220495044085
$ws.Range("T10").Value = 0.0001224220495044085
$ws.Range("G11").Value = 2.162051666666667
$ws.Range("H11").Value = 6.486155
$ws.Range("I11").Value = 0.03207188515291837
$ws.Range("J11").Value = 0.03207188515291837
$ws.Range("O11").Value = 0.9945745510447523
$ws.Range("P11").Value = 0.9945745510447522
$ws.Range("Q11").Value = 375.805471991206
$ws.Range("R11").Value = 3382.249247920855
$ws.Range("S11").Value = 0.03189788077712264
$ws.Range("T11").Value = 0.03189788077712264
$ws.Range("G12").Value = 2.162051666666667
$ws.Range("H12").Value = 6.486155
$ws.Range("I12").Value = 0.03207188515291837
$ws.Range("J12").Value = 0.03207188515291837
$ws.Range("K12").Value = 2.0
$ws.Range("L12").Value = 0.6666666666666666
$ws.Range("M12").Value = 0.260372
$ws.Range("N12").Value = 0.7811159999999999
$ws.Range("O12").Value = 0.001489822709979835
$ws.Range("P12").Value = 0.001489822709979834
$ws.Range("Q12").Value = 0.5629377165533334
$ws.Range("R12").Value = 5.06643944898
$ws.Range("S12").Value = 0.00004778142285268286
$ws.Range("T12").Value = 0.00004778142285268286
$ws.Range("G13").Value = 2.162051666666667
$ws.Range("H13").Value = 6.486155
$ws.Range("I13").Value = 0.03207188515291837
$ws.Range("J13").Value = 0.03207188515291837
$ws.Range("K13").Value = 2.0
$ws.Range("L13").Value = 0.6666666666666666
$ws.Range("M13").Value = 0.020712
$ws.Range("N13").Value = 0.062136
$ws.Range("O13").Value = 0.0001185120057805845
$ws.Range("P13").Value = 0.0001185120057805844
$ws.Range("Q13").Value = 0.04478041411999999
$ws.Range("R13").Value = 0.40302372708
$ws.Range("S13").Value = 0.000003800903438636903
$ws.Range("T13").Value = 0.000003800903438636902
$ws.Range("G14").Value = 13.28262333333333
$ws.Range("H14").Value = 39.84787
$ws.Range("I14").Value = 0.1970345004441647
$ws.Range("J14").Value = 0.1970345004441647
$ws.Range("M14").Value = 0.667106
$ws.Range("N14").Value = 2.001318
$ws.Range("O14").Value = 0.003817114239487378
$ws.Range("P14").Value = 0.003817114239487378
$ws.Range("Q14").Value = 8.860917721406667
$ws.Range("R14").Value = 79.74825949266
$ws.Range("S14").Value = 0.0007521031973157033
$ws.Range("T14").Value = 0.0007521031973157032
$ws.Range("G15").Value = 13.28262333333333
$ws.Range("H15").Value = 39.84787
$ws.Range("I15").Value = 0.1970345004441647
$ws.Range("J15").Value = 0.1970345004441647
$ws.Range("O15").Value = 0.9945745510447523
$ws.Range("P15").Value = 0.9945745510447522
$ws.Range("Q15").Value = 2308.771158443519
$ws.Range("R15").Value = 20778.94042599167
$ws.Range("S15").Value = 0.1959654998195822
$ws.Range("T15").Value = 0.1959654998195822
$ws.Range("G16").Value = 13.28262333333333
$ws.Range("H16").Value = 39.84787
$ws.Range("I16").Value = 0.1970345004441647
$ws.Range("J16").Value = 0.1970345004441647
$ws.Range("K16").Value = 2.0
$ws.Range("L16").Value = 0.6666666666666666
$ws.Range("M16").Value = 0.260372
$ws.Range("N16").Value = 0.7811159999999999
$ws.Range("O16").Value = 0.001489822709979835
$ws.Range("P16").Value = 0.001489822709979834
$ws.Range("Q16").Value = 3.458423202546667
$ws.Range("R16").Value = 31.12580882292
$ws.Range("S16").Value = 0.0002935464734112485
$ws.Range("T16").Value = 0.0002935464734112484
$ws.Range("G17").Value = 13.28262333333333
$ws.Range("H17").Value = 39.84787
$ws.Range("I17").Value = 0.1970345004441647
$ws.Range("J17").Value = 0.1970345004441647
$ws.Range("K17").Value = 2.0
$ws.Range("L17").Value = 0.6666666666666666
$ws.Range("M17").Value = 0.020712
$ws.Range("N17").Value = 0.062136
$ws.Range("O17").Value = 0.0001185120057805845
$ws.Range("P17").Value = 0.0001185120057805844
$ws.Range("Q17").Value = 0.27510969448
$ws.Range("R17").Value = 2.47598725032
$ws.Range("S17").Value = 0.00002335095385561342
$ws.Range("T17").Value = 0.00002335095385561342
